$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.214.76"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "2.093.74"
$ws.Range("E3").Value = "  +9.39%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("E6").Value = "  -3.82%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.99"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.374"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0749"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("E12").Value = "  +5.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("D14").Value = "2.396.67"
$ws.Range("E14").Value = "  +9.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.836"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").Value = "2.087.65"
$ws.Range("E16").Value = "  +8.65%  "
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "37.038.43"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.40%  "
$ws.Range("D20").Value = "0.0₃0825"
$ws.Range("E20").Value = "  -3.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "240.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.56%  "
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.46%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.45%  "
$ws.Range("E29").Value = "  -4.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +40.45%  "
$ws.Range("E31").Value = "  -4.41%  "
$ws.Range("E32").Value = "  +26.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0608"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0933"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.03%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("E37").Value = "  +18.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.60%  "
$ws.Range("E39").Value = "  -3.37%  "
$ws.Range("E40").Value = "  -8.40%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.41%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0223"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("E43").Value = "  +6.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.75%  "
$ws.Range("E45").Value = "  -2.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0870"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.59%  "
$ws.Range("B47").Value = "MXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.13%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.318.72"
$ws.Range("E48").Value = "  -2.14%  "
$ws.Range("E49").Value = "  +8.88%  "
$ws.Range("D50").Value = "2.270.45"
$ws.Range("E50").Value = "  +8.74%  "
$ws.Range("E51").Value = "  -5.39%  "
